$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.631.39"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "1.597.90"

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.65"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("E8").Value = "  -0.54%  "

$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.57"
$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").Value = "1.822.33"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").Value = "1.593.11"
$ws.Range("E13").Value = "  +0.31%  "

$ws.Range("E14").Value = "  +0.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.60"
$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").Value = "26.607.68"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("E18").Value = "  -2.50%  "

$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.27"

$ws.Range("E21").Value = "  -1.14%  "

$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.23"
$ws.Range("E23").Value = "  -3.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.74"
$ws.Range("E25").Value = "  +0.40%  "

$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("E27").Value = "  +0.32%  "

$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.28"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0507"
$ws.Range("E30").Value = "  -1.92%  "

$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("E34").Value = "  +18.68%  "

$ws.Range("D35").Value = "1.278.21"
$ws.Range("E35").Value = "  -0.91%  "

$ws.Range("E36").Value = "  +0.81%  "

$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.598"
$ws.Range("E38").Value = "  -3.26%  "

$ws.Range("E39").Value = "  -1.91%  "

$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.17"
$ws.Range("E42").Value = "  -0.98%  "

$ws.Range("E43").Value = "  -1.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.66"
$ws.Range("E44").Value = "  -0.75%  "

$ws.Range("D45").Value = "1.733.93"
$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.62"
$ws.Range("E46").Value = "  -1.57%  "

$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.103"
$ws.Range("E48").Value = "  +2.05%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.44"
$ws.Range("E51").Value = "  +1.20%  "
